# This script reproduces the daily cryptos-list refresh committed by the
# "Updated cryptos list ... with GitHub Actions" workflow: it overwrites the
# Price (column D) and Volume(1h) (column E) figures for the existing coin
# rows with the newly scraped figures, and fixes the ranking swap between
# LidoDAOToken and MXToken (rows 37/38, including their Coin name and Link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimal numbers (e.g. "0.998").
# Excel auto-detects such strings and silently turns the cell into a numeric
# cell (which also mangles the value through binary floating point, e.g.
# "43.58" -> 43.579999999999998). The source workbook stores every Price cell
# as text, so those particular cells are temporarily switched to the Text
# number format before the value is written, then switched back to the
# workbook default style (Normal) so no stray formatting is left behind.
$textForceRange = $excel.Union($ws.Range("D4"), $ws.Range("D5"), $ws.Range("D6"), $ws.Range("D8"), $ws.Range("D9"), $ws.Range("D11"), $ws.Range("D16"), $ws.Range("D17"), $ws.Range("D18"), $ws.Range("D19"), $ws.Range("D26"), $ws.Range("D27"), $ws.Range("D29"), $ws.Range("D31"), $ws.Range("D32"), $ws.Range("D35"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D39"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D45"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D50"))
foreach ($area in $textForceRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range('D2').Value = '29.459.16'
$ws.Range('E2').Value = '  +3.42%  '
$ws.Range('D3').Value = '1.604.12'
$ws.Range('E3').Value = '  +3.22%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '213.05'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '0.517'
$ws.Range('E6').Value = '  +6.85%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '26.66'
$ws.Range('E8').Value = '  +10.10%  '
$ws.Range('D9').Value = '43.58'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('D11').Value = '0.0596'
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').Value = '1.833.43'
$ws.Range('E13').Value = '  +3.24%  '
$ws.Range('D14').Value = '1.592.20'
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '29.494.55'
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('D16').Value = '0.536'
$ws.Range('E16').Value = '  +4.97%  '
$ws.Range('D17').Value = '3.74'
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').Value = '63.41'
$ws.Range('E18').Value = '  +3.52%  '
$ws.Range('D19').Value = '239.64'
$ws.Range('E19').Value = '  +4.41%  '
$ws.Range('E20').Value = '  +2.79%  '
$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  +3.14%  '
$ws.Range('E23').Value = '  +3.06%  '
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').Value = '154.43'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').Value = '15.25'
$ws.Range('E27').Value = '  +3.32%  '
$ws.Range('E28').Value = '  +4.40%  '
$ws.Range('D29').Value = '6.35'
$ws.Range('E29').Value = '  +1.73%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').Value = '0.0470'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '1.06'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('D34').Value = '1.429.92'
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('D35').Value = '3.09'
$ws.Range('E35').Value = '  +2.80%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '1.51'
$ws.Range('E37').Value = '  +1.71%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.80'
$ws.Range('E38').Value = '  +4.94%  '
$ws.Range('D39').Value = '2.28'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('E41').Value = '  +2.91%  '
$ws.Range('D42').Value = '1.93'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').Value = '53.41'
$ws.Range('E43').Value = '  +22.33%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '0.793'
$ws.Range('E45').Value = '  +2.54%  '
$ws.Range('D46').Value = '0.0473'
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('D47').Value = '65.72'
$ws.Range('E47').Value = '  +2.94%  '
$ws.Range('D48').Value = '5.26'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').Value = '1.744.20'
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('D50').Value = '86.68'
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('E51').Value = '  -3.63%  '

# Restore the default (Normal) style on the cells we temporarily reformatted
# above, now that their text values are safely stored.
foreach ($area in $textForceRange.Areas) {
    $area.Style = "Normal"
}
